$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ------------------------------------------------------------------
# 1) Remove the (empty) Sheet3
# ------------------------------------------------------------------
$wb.Worksheets("Sheet3").Delete()

# ------------------------------------------------------------------
# 2) Update font (Sheet1 / Sheet2 data range A1:D117) from the
#    Chinese "Song" font to "Times New Roman", and bump the row
#    height of the populated rows from 13.5 to 15 (the natural
#    Excel default height once the font changes away from the
#    East-Asian font), and widen the columns that needed more room.
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets("Sheet1")
$ws1.Range("A1:D117").Font.Name = "Times New Roman"
$ws1.Rows("1:117").RowHeight = 15
$ws1.Columns("C").ColumnWidth = 10.43

$ws2 = $wb.Worksheets("Sheet2")
$ws2.Range("A1:D117").Font.Name = "Times New Roman"
$ws2.Rows("1:117").RowHeight = 15
$ws2.Columns("C:D").ColumnWidth = 10.43

Write-Host "done"
